$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "Importar"
$ws.Range("G4").ClearContents()
$ws.Range("B7").Select()
